# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.326.88"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.932.19"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'0.7514"
$ws.Range("E5").Value = "  +5.28%  "
$ws.Range("D6").Value = "'244.89"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.3182"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("D9").Value = "'27.60"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'0.06983"
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("D11").Value = "'0.7804"
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("D12").Value = "'0.07988"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'1.927.28"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'5.349"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "'94.29"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "'14.39"
$ws.Range("E16").Value = "  -3.42%  "
$ws.Range("D17").Value = "'30.325.39"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'252.75"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'0.000007929"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").Value = "'5.720"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "'2.189.52"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'6.671"
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("D25").Value = "'9.478"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "'165.95"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'18.93"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "'0.1327"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").Value = "'2.218"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").Value = "'1.361"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "'1.511"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").Value = "'4.372"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").Value = "'4.111"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").Value = "'0.05156"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'1.271"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'0.7442"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'0.01944"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "'2.793"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'77.81"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "'6.396"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'0.4459"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "'1.964"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'0.8315"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "'100.78"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "'9.741"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'7.459"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "'985.12"
$ws.Range("E49").Value = "  +11.78%  "
$ws.Range("D50").Value = "'37.13"
$ws.Range("D51").Value = "'0.06005"
$ws.Range("E51").Value = "  -0.86%  "
